# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 207 (shifting the existing
# rows 207-251 down to 208-252) in the Feria Lagunitas de Puerto Montt - Uva
# price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207; this shifts existing rows 207..251
# down to 208..252 and keeps their data/formatting intact.
$ws.Rows("207").Insert()

# Populate the newly inserted row 207 with the new price record.
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").Value = 44782
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100109
$ws.Range("H207").Value = "Uva"
$ws.Range("I207").Value = 100109001
$ws.Range("J207").Value = "Uva"
$ws.Range("K207").Value = "Red Globe"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 300
$ws.Range("N207").Value = 15000
$ws.Range("O207").Value = 16000
$ws.Range("P207").Value = 15500
$ws.Range("Q207").Value = "$/bandeja 8 kilos"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 1938
$ws.Range("T207").Value = 8
